$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (values only change)
$ws.Range("B2").Value = -6.376377285340954
$ws.Range("C2").Value = -6.376377285340955
$ws.Range("D2").Value = -6.376377285340955

# Row 3 - RandomForestRegressor (values only change)
$ws.Range("B3").Value = 0.845049961883345
$ws.Range("C3").Value = 0.9351927071930313
$ws.Range("D3").Value = 0.6502301898721462

# Row 4 - renamed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9373413926863283
$ws.Range("C4").Value = 0.9108926000615768
$ws.Range("D4").Value = 0.7936880044057121

# Row 5 - renamed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8161247051948194
$ws.Range("C5").Value = 0.396937205572043
$ws.Range("D5").Value = -3.498763047294248
